# Apply cryptos list update (price/volume refresh + two coin-row swaps)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value2 = "36.578.62"
$ws.Cells.Item(2, 5).Value2 = "  -0.66%  "

$ws.Cells.Item(3, 4).Value2 = "2.091.31"
$ws.Cells.Item(3, 5).Value2 = "  +9.22%  "

$ws.Cells.Item(4, 5).Value2 = "  +0.12%  "

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value2 = "252.14"
$ws.Cells.Item(5, 5).Value2 = "  +0.63%  "

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value2 = "0.659"
$ws.Cells.Item(6, 5).Value2 = "  -6.11%  "

$ws.Cells.Item(7, 5).Value2 = "  +0.11%  "

$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value2 = "48.99"
$ws.Cells.Item(8, 5).Value2 = "  +4.53%  "

$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value2 = "60.22"
$ws.Cells.Item(9, 5).Value2 = "  +3.64%  "

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value2 = "0.375"
$ws.Cells.Item(10, 5).Value2 = "  +0.44%  "

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value2 = "0.0745"
$ws.Cells.Item(11, 5).Value2 = "  -2.15%  "

$ws.Cells.Item(12, 5).Value2 = "  +0.16%  "

$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value2 = "14.76"
$ws.Cells.Item(13, 5).Value2 = "  +0.30%  "

$ws.Cells.Item(14, 4).Value2 = "2.405.91"
$ws.Cells.Item(14, 5).Value2 = "  +9.76%  "

$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value2 = "0.833"
$ws.Cells.Item(15, 5).Value2 = "  +1.78%  "

$ws.Cells.Item(16, 4).Value2 = "2.099.30"
$ws.Cells.Item(16, 5).Value2 = "  +9.60%  "

$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value2 = "5.11"
$ws.Cells.Item(17, 5).Value2 = "  -0.47%  "

$ws.Cells.Item(18, 4).Value2 = "36.545.94"
$ws.Cells.Item(18, 5).Value2 = "  -1.83%  "

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value2 = "72.96"
$ws.Cells.Item(19, 5).Value2 = "  -2.42%  "

$ws.Cells.Item(20, 4).Value2 = "0.0₃0834"
$ws.Cells.Item(20, 5).Value2 = "  -2.98%  "

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value2 = "13.25"
$ws.Cells.Item(21, 5).Value2 = "  -2.82%  "

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value2 = "239.89"
$ws.Cells.Item(22, 5).Value2 = "  -4.61%  "

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value2 = "5.28"
$ws.Cells.Item(23, 5).Value2 = "  +1.58%  "

$ws.Cells.Item(24, 5).Value2 = "  -0.11%  "

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value2 = "2.53"
$ws.Cells.Item(25, 5).Value2 = "  -3.44%  "

$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value2 = "170.93"
$ws.Cells.Item(26, 5).Value2 = "  +1.72%  "

$ws.Cells.Item(27, 2).Value2 = "Cosmos"
$ws.Cells.Item(27, 3).Value2 = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value2 = "9.26"
$ws.Cells.Item(27, 5).Value2 = "  +4.98%  "

$ws.Cells.Item(28, 2).Value2 = "EthereumClassic"
$ws.Cells.Item(28, 3).Value2 = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value2 = "21.13"
$ws.Cells.Item(28, 5).Value2 = "  +12.86%  "

$ws.Cells.Item(29, 5).Value2 = "  -10.26%  "

$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value2 = "28.42"
$ws.Cells.Item(30, 5).Value2 = "  +48.34%  "

$ws.Cells.Item(31, 5).Value2 = "  -5.40%  "

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value2 = "4.49"
$ws.Cells.Item(32, 5).Value2 = "  -2.88%  "

$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value2 = "0.0617"
$ws.Cells.Item(33, 5).Value2 = "  -0.65%  "

$ws.Cells.Item(34, 2).Value2 = "ImmutableX"
$ws.Cells.Item(34, 3).Value2 = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value2 = "0.980"
$ws.Cells.Item(34, 5).Value2 = "  +11.75%  "

$ws.Cells.Item(35, 2).Value2 = "LidoDAOToken"
$ws.Cells.Item(35, 3).Value2 = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value2 = "2.41"
$ws.Cells.Item(35, 5).Value2 = "  +19.45%  "

$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value2 = "0.0901"
$ws.Cells.Item(36, 5).Value2 = "  +0.50%  "

$ws.Cells.Item(37, 5).Value2 = "  +0.16%  "

$ws.Cells.Item(38, 5).Value2 = "  -2.97%  "

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value2 = "4.09"
$ws.Cells.Item(39, 5).Value2 = "  -5.99%  "

$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value2 = "1.35"
$ws.Cells.Item(40, 5).Value2 = "  -10.93%  "

$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value2 = "0.0224"
$ws.Cells.Item(41, 5).Value2 = "  -1.91%  "

$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value2 = "1.16"
$ws.Cells.Item(42, 5).Value2 = "  +5.47%  "

$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value2 = "97.71"
$ws.Cells.Item(43, 5).Value2 = "  -7.54%  "

$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value2 = "16.46"
$ws.Cells.Item(44, 5).Value2 = "  -8.61%  "

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value2 = "2.79"
$ws.Cells.Item(45, 5).Value2 = "  -3.14%  "

$ws.Cells.Item(46, 4).Value2 = "1.335.40"
$ws.Cells.Item(46, 5).Value2 = "  -1.00%  "

$ws.Cells.Item(47, 5).Value2 = "  +3.46%  "

$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value2 = "7.09"
$ws.Cells.Item(48, 5).Value2 = "  +9.31%  "

$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value2 = "2.87"
$ws.Cells.Item(49, 5).Value2 = "  +1.38%  "

$ws.Cells.Item(50, 4).Value2 = "2.303.00"
$ws.Cells.Item(50, 5).Value2 = "  +10.21%  "

$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value2 = "2.25"
$ws.Cells.Item(51, 5).Value2 = "  -6.55%  "
